# "added harvard case classification" -- refresh the Avey/WebMD_old stats table:
#   - header row 1: swap the average_doctor / average_doctor_old column labels
#     (BP/BQ) so they line up with the recomputed data beneath them
#   - rows 4-13: update the recomputed precision/recall/F-score/NDCG/M1/M3/M5/
#     position/length stats that shifted once the Harvard case classification
#     was folded into the "K health_old", "WebMD_old" and "*_old" doctor columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header labels (row 1): average_doctor / average_doctor_old were swapped ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- row 4: stats_for_precision ---
$ws.Range("E4").Value = 0.429
$ws.Range("F4").Value = 0.071
$ws.Range("G4").Value = 0.267
$ws.Range("N4").Value = 0.438
$ws.Range("O4").Value = 0.065
$ws.Range("P4").Value = 0.255
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.131
$ws.Range("W4").Value = 0.292
$ws.Range("X4").Value = 0.111
$ws.Range("Y4").Value = 0.333
$ws.Range("AI4").Value = 0.288
$ws.Range("AJ4").Value = 0.088
$ws.Range("AK4").Value = 0.297
$ws.Range("AU4").Value = 0.196
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.166
$ws.Range("BA4").Value = 1.973
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.398
$ws.Range("BG4").Value = 0.72
$ws.Range("BH4").Value = 0.144
$ws.Range("BI4").Value = 0.379
$ws.Range("BM4").Value = 0.706
$ws.Range("BN4").Value = 0.08
$ws.Range("BO4").Value = 0.282
$ws.Range("BP4").Value = 0.658
$ws.Range("BQ4").Value = 0.701

# --- row 5: stats_for_recall ---
$ws.Range("E5").Value = 0.544
$ws.Range("F5").Value = 0.086
$ws.Range("G5").Value = 0.293
$ws.Range("N5").Value = 0.741
$ws.Range("O5").Value = 0.077
$ws.Range("P5").Value = 0.278
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.084
$ws.Range("W5").Value = 0.279
$ws.Range("X5").Value = 0.11
$ws.Range("Y5").Value = 0.332
$ws.Range("AI5").Value = 0.306
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.308
$ws.Range("AU5").Value = 0.381
$ws.Range("AV5").Value = 0.091
$ws.Range("AW5").Value = 0.302
$ws.Range("BA5").Value = 1.331
$ws.Range("BB5").Value = 0.082
$ws.Range("BC5").Value = 0.286
$ws.Range("BG5").Value = 0.391
$ws.Range("BH5").Value = 0.052
$ws.Range("BI5").Value = 0.228
$ws.Range("BM5").Value = 0.553
$ws.Range("BN5").Value = 0.066
$ws.Range("BO5").Value = 0.256
$ws.Range("BP5").Value = 0.444
$ws.Range("BQ5").Value = 0.454

# --- row 6: stats_for_f1-score ---
$ws.Range("E6").Value = 0.48
$ws.Range("N6").Value = 0.551
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.285
$ws.Range("AI6").Value = 0.297
$ws.Range("AU6").Value = 0.259
$ws.Range("BA6").Value = 1.58
$ws.Range("BG6").Value = 0.507
$ws.Range("BM6").Value = 0.62
$ws.Range("BP6").Value = 0.527
$ws.Range("BQ6").Value = 0.548

# --- row 7: stats_for_f2-score ---
$ws.Range("E7").Value = 0.516
$ws.Range("N7").Value = 0.651
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.282
$ws.Range("AI7").Value = 0.302
$ws.Range("AU7").Value = 0.32
$ws.Range("BA7").Value = 1.419
$ws.Range("BG7").Value = 0.43
$ws.Range("BM7").Value = 0.578
$ws.Range("BP7").Value = 0.473
$ws.Range("BQ7").Value = 0.487

# --- row 8: stats_for_NDCG ---
$ws.Range("E8").Value = 0.604
$ws.Range("F8").Value = 0.111
$ws.Range("G8").Value = 0.333
$ws.Range("N8").Value = 0.778
$ws.Range("O8").Value = 0.061
$ws.Range("P8").Value = 0.248
$ws.Range("Q8").Value = 0.018
$ws.Range("W8").Value = 0.305
$ws.Range("AI8").Value = 0.325
$ws.Range("AJ8").Value = 0.125
$ws.Range("AK8").Value = 0.353
$ws.Range("AU8").Value = 0.322
$ws.Range("AW8").Value = 0.29
$ws.Range("BA8").Value = 1.724
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.353
$ws.Range("BG8").Value = 0.556
$ws.Range("BH8").Value = 0.108
$ws.Range("BI8").Value = 0.328
$ws.Range("BM8").Value = 0.691
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.259
$ws.Range("BP8").Value = 0.575
$ws.Range("BQ8").Value = 0.599

# --- row 9: stats_for_M1 ---
$ws.Range("E9").Value = 0.544
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.678
$ws.Range("O9").Value = 0.218
$ws.Range("P9").Value = 0.467
$ws.Range("W9").Value = 0.2
$ws.Range("X9").Value = 0.16
$ws.Range("Y9").Value = 0.4
$ws.Range("AI9").Value = 0.244
$ws.Range("AJ9").Value = 0.185
$ws.Range("AK9").Value = 0.43
$ws.Range("BA9").Value = 1.666
$ws.Range("BB9").Value = 0.246
$ws.Range("BC9").Value = 0.496
$ws.Range("BG9").Value = 0.589
$ws.Range("BH9").Value = 0.242
$ws.Range("BI9").Value = 0.492
$ws.Range("BM9").Value = 0.644
$ws.Range("BN9").Value = 0.229
$ws.Range("BO9").Value = 0.479
$ws.Range("BP9").Value = 0.555
$ws.Range("BQ9").Value = 0.577

# --- row 10: stats_for_M3 ---
$ws.Range("E10").Value = 0.678
$ws.Range("F10").Value = 0.218
$ws.Range("G10").Value = 0.467
$ws.Range("N10").Value = 0.878
$ws.Range("O10").Value = 0.107
$ws.Range("P10").Value = 0.328
$ws.Range("W10").Value = 0.378
$ws.Range("X10").Value = 0.235
$ws.Range("Y10").Value = 0.485
$ws.Range("AI10").Value = 0.356
$ws.Range("AJ10").Value = 0.229
$ws.Range("AK10").Value = 0.479
$ws.Range("AU10").Value = 0.311
$ws.Range("AV10").Value = 0.214
$ws.Range("AW10").Value = 0.463
$ws.Range("BA10").Value = 2.055
$ws.Range("BB10").Value = 0.246
$ws.Range("BC10").Value = 0.496
$ws.Range("BG10").Value = 0.644
$ws.Range("BH10").Value = 0.229
$ws.Range("BI10").Value = 0.479
$ws.Range("BM10").Value = 0.844
$ws.Range("BN10").Value = 0.131
$ws.Range("BO10").Value = 0.362
$ws.Range("BP10").Value = 0.685
$ws.Range("BQ10").Value = 0.719

# --- row 11: stats_for_M5 ---
$ws.Range("E11").Value = 0.711
$ws.Range("F11").Value = 0.205
$ws.Range("G11").Value = 0.453
$ws.Range("N11").Value = 0.9
$ws.Range("O11").Value = 0.09
$ws.Range("P11").Value = 0.3
$ws.Range("W11").Value = 0.378
$ws.Range("X11").Value = 0.235
$ws.Range("Y11").Value = 0.485
$ws.Range("AI11").Value = 0.389
$ws.Range("AJ11").Value = 0.238
$ws.Range("AK11").Value = 0.487
$ws.Range("AU11").Value = 0.456
$ws.Range("AV11").Value = 0.248
$ws.Range("AW11").Value = 0.498
$ws.Range("BA11").Value = 2.055
$ws.Range("BB11").Value = 0.246
$ws.Range("BC11").Value = 0.496
$ws.Range("BG11").Value = 0.644
$ws.Range("BH11").Value = 0.229
$ws.Range("BI11").Value = 0.479
$ws.Range("BM11").Value = 0.844
$ws.Range("BN11").Value = 0.131
$ws.Range("BO11").Value = 0.362
$ws.Range("BP11").Value = 0.685
$ws.Range("BQ11").Value = 0.722

# --- row 12: stats_for_position ---
$ws.Range("E12").Value = 1.422
$ws.Range("F12").Value = 0.775
$ws.Range("G12").Value = 0.88
$ws.Range("N12").Value = 1.482
$ws.Range("O12").Value = 1.069
$ws.Range("P12").Value = 1.034
$ws.Range("W12").Value = 1.647
$ws.Range("X12").Value = 0.581
$ws.Range("Y12").Value = 0.762
$ws.Range("AI12").Value = 1.743
$ws.Range("AJ12").Value = 1.334
$ws.Range("AK12").Value = 1.155
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.739
$ws.Range("BB12").Value = 0.418
$ws.Range("BC12").Value = 0.647
$ws.Range("BG12").Value = 1.103
$ws.Range("BH12").Value = 0.127
$ws.Range("BI12").Value = 0.357
$ws.Range("BM12").Value = 1.303
$ws.Range("BN12").Value = 0.343
$ws.Range("BO12").Value = 0.585
$ws.Range("BP12").Value = 1.246
$ws.Range("BQ12").Value = 1.269

# --- row 13: stats_for_length (x of gs) ---
$ws.Range("E13").Value = 1.57
$ws.Range("F13").Value = 0.662
$ws.Range("G13").Value = 0.814
$ws.Range("N13").Value = 2.071
$ws.Range("O13").Value = 0.964
$ws.Range("P13").Value = 0.982
$ws.Range("W13").Value = 1.025
$ws.Range("X13").Value = 0.192
$ws.Range("Y13").Value = 0.438
$ws.Range("AI13").Value = 1.281
$ws.Range("AJ13").Value = 0.378
$ws.Range("AK13").Value = 0.615
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.38
$ws.Range("BB13").Value = 0.304
$ws.Range("BC13").Value = 0.551
$ws.Range("BG13").Value = 0.586
$ws.Range("BH13").Value = 0.073
$ws.Range("BI13").Value = 0.271
$ws.Range("BM13").Value = 0.911
$ws.Range("BN13").Value = 0.285
$ws.Range("BO13").Value = 0.534
$ws.Range("BP13").Value = 0.793
$ws.Range("BQ13").Value = 0.732
